# Changing how capacity market is considered in investment module
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coupling Parameters")

# End Year: 2090 -> 2055
$ws.Range("B4").Value = 2055

# start_dismantling_tick: 1 -> 5
$ws.Range("B23").Value = 5

# capacity_remuneration_mechanism: strategic_reserve -> capacity_market
$ws.Range("B44").Value = "capacity_market"

# Update the view to match the new scroll/selection position
$ws.Application.Goto($ws.Range("A13"), $false)
$ws.Range("C34").Select()
